$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A371").Value = '#Kenya Due to massive impact of #Covid19 on food systems and food security, conversations on food issues that have always been left out; eg farm insurance, nutrition based agriculture, social protection, food access, self reliance; should happen. #Covid19FoodFuture @sharonjcheboi'
$ws.Range("A373").Value = '#Kenya Post Covid19, countries should invest in strong identification systems for easier and efficient food donation programmes in the future. #Covid19FoodFuture @sharonjcheboi'
$ws.Range("A375").Value = '#Kenya Food pattern change is a short term impact of #Covid19 on consumption behaviours. Kenyans have begun to prefer locally produced fish to imports. Read more on impacts on small scale traders in #Nairobi by @jamessmat https://bit.ly/2Y3U9PI  @sharonjcheboi #Covid19FoodFuture'
$ws.Range("A377").Value = '#Kenya should never have to choose between death by #Covid19 or hunger. Challenges #SouthAfrica has faced in food aid are similar to those in Kenya. This should be used to develop solutions to anticipated challenges. Read: https://bit.ly/3bqK5Eq @sharonjcheboi #Covid19FoodFuture'
$ws.Range("A379").Value = '#Kenya adoption of digital solutions to feed the most vulnerable during #Covid19 is necessary to adhere to government guidelines on social distancing. Yet, it remains underutilized. Read more: https://bit.ly/2VSWJFN by @SerahKiragu_tmg #Covid19FoodFuture  @sharonjcheboi'
$ws.Range("A381").Value = '#Kenya imports more food than it exports. Global food prices have increased significantly and a high population is made up of low income persons. With lack of proper social protection, chronic hunger might occur due to low purchasing power. #covid19foodfuture @sharonjcheboi'
$ws.Range("A383").Value = '#Kenya food markets have been destabilized as expressed in food shortages and massive post harvest losses. Poor market linkage is now a big challenge . Linking producers to consumers is necessary to stabilize food markets and reduce uncertainty. #Covid19FoodFuture @sharonjcheboi'
$ws.Range("A385").Value = 'Staple food production in developed countries have been affected less by #Covid19. #Africa largely consists of small scale farmers dependent on human labour.There''s more need to strengthen small scale farmers. Here''s how: https://bit.ly/3bli4OC #Covid19FoodFuture @sharonjcheboi'
$ws.Range("A387").Value = 'While #Africa already suffers from food insecurity, it cannot afford a food crisis amidst a #Covid19 health crisis. Locking down farmers from food production will lead to food shortage for the increasing population. #Covid19FoodFuture @sharonjcheboi'
$ws.Range("A389").Value = '#Kenya government focus has shifted to cushion the health system through: Masks, Employing medics, Medical units Less attention has been paid to #Covid19 driven impacts on food security. This is a very risky rope to thread on. #Covid19FoodFuture @sharonjcheboi'
$ws.Range("A391").Value = 'In #Africa, food supply chains and logistics have been heavily affected as a result of: closed borders, national lockdowns disrupted transport system. The impact is being felt by both producers and consumers. #Covid19FoodFuture @sharonjcheboi'
$ws.Range("A393").Value = 'The impact of #Covid19 has been felt by the lowest income households first. Through higher food costs, less purchasing power, poor ability to stock due to income loss and inadequate refrigeration and lack of savings.  @sharonjcheboi #CovidFoodFuture'
$ws.Range("A395").Value = 'What is the extent of our nutrition sensitive agriculture? And do care givers have knowledge and socioeconomic capacity to ensure that? What harms has ultra-processed foods caused to health amidst #Covid_19? @sharonjcheboi #CovidFoodFuture'

$ws.Range("B370").Select()

